$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Issue number, reporting week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Plain numeric edits (value only, style/type unchanged) ---
$ws.Range("H14").Value = 100
$ws.Range("K14").Value = 100
$ws.Range("N14").Value = 100
$ws.Range("F15").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 8
$ws.Range("K16").Value = 60
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -85.185185185185
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = -30.76923076923
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 12.5
$ws.Range("N17").Value = -30.76923076923
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 77.777777777777
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 120
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -68.571428571428
$ws.Range("N18").Value = -90
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 4.444444444444
$ws.Range("I19").Value = 35
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 20.689655172413
$ws.Range("M19").Value = 66.666666666666
$ws.Range("N19").Value = -22.222222222222
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 150
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 166.666666666667
$ws.Range("L20").Value = 300
$ws.Range("M20").Value = -11.111111111111
$ws.Range("N20").Value = -93.22033898305
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -11.538461538461
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 24.358974358974
$ws.Range("I21").Value = 75
$ws.Range("J21").Value = 62
$ws.Range("K21").Value = 20.967741935483
$ws.Range("L21").Value = 33.928571428571
$ws.Range("M21").Value = -16.666666666666
$ws.Range("N21").Value = -78.005865102639
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -7.692307692307
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -33.734939759036
$ws.Range("I24").Value = 41
$ws.Range("J24").Value = 67
$ws.Range("K24").Value = -38.805970149253
$ws.Range("L24").Value = -18
$ws.Range("M24").Value = -38.805970149253
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 8.333333333333
$ws.Range("I25").Value = 19
$ws.Range("J25").Value = 15
$ws.Range("K25").Value = 26.666666666666
$ws.Range("L25").Value = 18.75
$ws.Range("M25").Value = -26.923076923076
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 700
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0

# --- Text -> Number conversions (copy a same-kind numeric cell's format first, so the
#     existing #,##0 / #,##0.0 style id is reused instead of a new style being created) ---
$intFormatSource = "G14"   # style 15 (#,##0)
$decFormatSource = "H14"   # style 16 (#,##0.0;"-"#,##0.0)
$ws.Range($intFormatSource).Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 2
$ws.Range($intFormatSource).Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 2
$ws.Range($intFormatSource).Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 2
$ws.Range($decFormatSource).Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("M15").Value = 100
$ws.Range($intFormatSource).Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range($decFormatSource).Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range($intFormatSource).Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J27").Value = 1
$ws.Range($decFormatSource).Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 700

# --- Number -> Text conversions ("0" / "***.*" placeholders) ---
# A leading apostrophe forces literal text instead of Excel auto-coercing the numeric-
# looking string back to a number; re-pasting the format from a known text cell (style 14)
# afterwards clears the transient quote-prefix format Excel applies on the literal assign.
$textFormatSource = "D26"   # style 14 (General, text)
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("C26").Value = "'0"
$ws.Range("C27").Value = "'0"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").Value = "'***.*"
$ws.Range("G29").Value = "'0"
$ws.Range("H29").Value = "'***.*"
$ws.Range("C30").Value = "'0"
$ws.Range($textFormatSource).Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range($textFormatSource).Copy()
$ws.Range("C30").PasteSpecial(-4122)
